$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2024-03-22 Friday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-03-23 Saturday", 2)

# Update the division problems in the table, row by row / cell by cell to
# avoid ambiguity since several problems share identical text (e.g. "17÷6=",
# "88÷5="). Rows 1, 5, 9, 13, 17 (1-based) contain the problems; the rows in
# between are blank answer rows.

$tbl = $d.Tables.Item(1)

$rows = @(1, 5, 9, 13, 17)
$values = @(
    @("17÷6=", "40÷9=", "49÷8=", "30÷4=", "61÷5="),
    @("75÷4=", "99÷5=", "11÷6=", "72÷4=", "13÷4="),
    @("32÷2=", "99÷9=", "28÷2=", "53÷5=", "20÷7="),
    @("70÷5=", "45÷3=", "73÷4=", "96÷9=", "37÷5="),
    @("39÷6=", "64÷6=", "84÷4=", "55÷6=", "99÷6=")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    for ($c = 1; $c -le 5; $c++) {
        $cell = $tbl.Cell($r, $c)
        $rng = $cell.Range
        $rng.End = $rng.End - 1
        $rng.Text = $values[$i][$c - 1]
    }
}
